# Insert a new data row above row 185 (pushing the existing rows 185-301
# down to 186-302) and populate the new row with the latest price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(185).Insert()

$ws.Range("A185").Value = 10
$ws.Range("B185").Value = "Vega Modelo de Temuco"
$ws.Range("C185").Value = "La Araucanía"
$ws.Range("D185").Value = 45086
$ws.Range("E185").Value = 9
$ws.Range("F185").Value = 100112005
$ws.Range("G185").Value = "Puerro"
$ws.Range("H185").Value = "Azul de Maquehue"
$ws.Range("I185").Value = "Primera"
$ws.Range("J185").Value = 30
$ws.Range("K185").Value = 10000
$ws.Range("L185").Value = 10000
$ws.Range("M185").Value = 10000
$ws.Range("N185").Value = "`$/docena de paquetes"
$ws.Range("O185").Value = "Provincia de Cautín"
$ws.Range("P185").Value = 833
$ws.Range("Q185").Value = 12
$ws.Range("R185").Value = "Hortaliza"
